$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported. It lands between the existing
# records at row 187 and 188, so insert a fresh row at position 188 and
# push every following record (old rows 188-215) down by one (new rows
# 189-216).
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new record's data.
$ws.Cells.Item(188, 1).Value = 3
$ws.Cells.Item(188, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 44491
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 100112009
$ws.Cells.Item(188, 7).Value = "Acelga"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 280
$ws.Cells.Item(188, 11).Value = 1800
$ws.Cells.Item(188, 12).Value = 2000
$ws.Cells.Item(188, 13).Value = 1914
$ws.Cells.Item(188, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(188, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(188, 16).Value = 319
$ws.Cells.Item(188, 17).Value = 6
$ws.Cells.Item(188, 18).Value = "Hortaliza"
